$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2919.3333
$ws.Range("I70").Value = 2929
$ws.Range("K70").Value = 8787
$ws.Range("M70").Value = -8517

# Row 73
$ws.Range("H73").Value = 2919.3333
$ws.Range("I73").Value = 2929
$ws.Range("K73").Value = 8787
$ws.Range("M73").Value = -7851

# Row 98
$ws.Range("H98").Value = 2182.04
$ws.Range("I98").Value = 1517.3636
$ws.Range("K98").Value = 1517.3636
$ws.Range("M98").Value = -19.36359999999991

# Row 100
$ws.Range("H100").Value = 1218.174
$ws.Range("I100").Value = 1218.174
$ws.Range("K100").Value = 1218.174
$ws.Range("M100").Value = -677.174

# Row 107
$ws.Range("H107").Value = 600.129
$ws.Range("I107").Value = 640.4815
$ws.Range("J107").Value = 327.75
$ws.Range("K107").Value = 640.4815
$ws.Range("L107").Value = 327.75
$ws.Range("M107").Value = 1279.5185
$ws.Range("N107").Value = -4167.75

# Row 116
$ws.Range("H116").Value = 5994.25
$ws.Range("I116").Value = 5992.3335
$ws.Range("K116").Value = 5992.3335
$ws.Range("M116").Value = -2550.3335

# Row 122
$ws.Range("H122").Value = 2182.04
$ws.Range("I122").Value = 1517.3636
$ws.Range("K122").Value = 4552.0908
$ws.Range("M122").Value = -2102.0908

# Row 132
$ws.Range("H132").Value = 1745.1409
$ws.Range("I132").Value = 1722.7761
$ws.Range("J132").Value = 2119.75
$ws.Range("K132").Value = 5168.3283
$ws.Range("L132").Value = 6359.25
$ws.Range("M132").Value = -2638.3283
$ws.Range("N132").Value = -11419.25

# Row 137
$ws.Range("H137").Value = 2451.158
$ws.Range("I137").Value = 2410.318
$ws.Range("J137").Value = 2507.3125
$ws.Range("K137").Value = 7230.954000000001
$ws.Range("L137").Value = 7521.9375
$ws.Range("M137").Value = -4680.954000000001
$ws.Range("N137").Value = -12621.9375

# Row 138
$ws.Range("H138").Value = 3944.4695
$ws.Range("I138").Value = 1969.6666
$ws.Range("J138").Value = 4584.946
$ws.Range("K138").Value = 5908.9998
$ws.Range("L138").Value = 13754.838
$ws.Range("M138").Value = -768.9997999999996
$ws.Range("N138").Value = -24034.838

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885

# Row 32
$ws.Range("H32").Value = 13936.32
$ws.Range("I32").Value = 7468.629
$ws.Range("J32").Value = 44782.23
$ws.Range("K32").Value = 7468.629
$ws.Range("L32").Value = 44782.23
$ws.Range("M32").Value = -7181.629
$ws.Range("N32").Value = -45356.23

# Row 45
$ws.Range("H45").Value = 11979.4
$ws.Range("I45").Value = 21669.4
$ws.Range("J45").Value = 2289.4
$ws.Range("K45").Value = 21669.4
$ws.Range("L45").Value = 2289.4
$ws.Range("M45").Value = -21292.4
$ws.Range("N45").Value = -3043.4

# Row 102
$ws.Range("H102").Value = 1142.55
$ws.Range("I102").Value = 1142.55
$ws.Range("K102").Value = 1142.55
$ws.Range("M102").Value = 479.45

# Row 110
$ws.Range("H110").Value = 5220.5938
$ws.Range("I110").Value = 5014.44
$ws.Range("J110").Value = 5956.857
$ws.Range("K110").Value = 5014.44
$ws.Range("L110").Value = 5956.857
$ws.Range("M110").Value = -2969.44
$ws.Range("N110").Value = -10046.857

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3127.5908
$ws.Range("I20").Value = 3021.6
$ws.Range("K20").Value = 3021.6
$ws.Range("M20").Value = -2774.6

# Row 107
$ws.Range("H107").Value = 3062.5557
$ws.Range("I107").Value = 3062.5557
$ws.Range("K107").Value = 3062.5557
$ws.Range("M107").Value = -1142.5557

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 1224.1818
$ws.Range("I7").Value = 1293.25
$ws.Range("K7").Value = 1293.25
$ws.Range("M7").Value = -1180.25

# Row 22
$ws.Range("H22").Value = 329.16666
$ws.Range("J22").Value = 331
$ws.Range("L22").Value = 331
$ws.Range("N22").Value = -1031

# Row 62
$ws.Range("H62").Value = 15296.333
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

# Row 65
$ws.Range("H65").Value = 15296.333
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

# Row 107
$ws.Range("H107").Value = 789.7368
$ws.Range("I107").Value = 699.6875
$ws.Range("K107").Value = 699.6875
$ws.Range("M107").Value = 1220.3125

$ws = $wb.Worksheets.Item("GSM")
# Row 17
$ws.Range("H17").Value = 320
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -836

# Row 52
$ws.Range("H52").Value = 32583.166
$ws.Range("I52").Value = 22749.5
$ws.Range("K52").Value = 22749.5
$ws.Range("M52").Value = -22490.5

# Row 126
$ws.Range("H126").Value = 9517.849
$ws.Range("I126").Value = 12173.608
$ws.Range("J126").Value = 3409.6
$ws.Range("K126").Value = 36520.824
$ws.Range("L126").Value = 10228.8
$ws.Range("M126").Value = -34050.824
$ws.Range("N126").Value = -15168.8

# Row 132
$ws.Range("H132").Value = 3247.6191
$ws.Range("I132").Value = 2473.818
$ws.Range("J132").Value = 4098.8
$ws.Range("K132").Value = 7421.454000000001
$ws.Range("L132").Value = 12296.4
$ws.Range("M132").Value = -4891.454000000001
$ws.Range("N132").Value = -17356.4

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4664.56
$ws.Range("I16").Value = 4400.6816
$ws.Range("K16").Value = 4400.6816
$ws.Range("M16").Value = -4230.6816

# Row 22
$ws.Range("H22").Value = 3196.074
$ws.Range("I22").Value = 2424.8333
$ws.Range("K22").Value = 2424.8333
$ws.Range("M22").Value = -2129.8333

# Row 27
$ws.Range("H27").Value = 3196.074
$ws.Range("I27").Value = 2424.8333
$ws.Range("K27").Value = 2424.8333
$ws.Range("M27").Value = -2317.8333

# Row 45
$ws.Range("H45").Value = 29563.334
$ws.Range("I45").Value = 29500
$ws.Range("K45").Value = 29500
$ws.Range("M45").Value = -29093

# Row 122
$ws.Range("H122").Value = 7784.3125
$ws.Range("I122").Value = 7600.636
$ws.Range("K122").Value = 22801.908
$ws.Range("M122").Value = -20351.908

# Row 132
$ws.Range("H132").Value = 5037.8945
$ws.Range("I132").Value = 4748.8335
$ws.Range("J132").Value = 5533.4287
$ws.Range("K132").Value = 14246.5005
$ws.Range("L132").Value = 16600.2861
$ws.Range("M132").Value = -11716.5005
$ws.Range("N132").Value = -21660.2861

# Row 141
$ws.Range("H141").Value = 98329.336
$ws.Range("J141").Value = 98329.336
$ws.Range("L141").Value = 98329.336
$ws.Range("N141").Value = -108689.336

$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 126
$ws.Range("H126").Value = 2508.8333
$ws.Range("I126").Value = 2157
$ws.Range("K126").Value = 6471
$ws.Range("M126").Value = -4001

# Row 141
$ws.Range("H141").Value = 92179.14
$ws.Range("J141").Value = 102920.8
$ws.Range("L141").Value = 102920.8
$ws.Range("N141").Value = -113280.8
